# Resolves #22 - rename header labels in the shared strings table.
# Only the text of the header cells (row 1) on every worksheet needs to
# change; the cells that reference them stay exactly the same.
#
# Mapping of old -> new header text:
#   sample.one      -> sample one
#   sample.one.1    -> sample one.1
#   sample.two      -> sample two
#   NA.             -> NA
#   sample.three    -> sample three
#   sample.three.1  -> sample three.1
#   X9..X12, X3..X8 -> 9..12, 3..8   (the leading "X" is dropped)

$map = @{
    "sample.one"     = "sample one"
    "sample.one.1"   = "sample one.1"
    "sample.two"     = "sample two"
    "NA."            = "NA"
    "sample.three"   = "sample three"
    "sample.three.1" = "sample three.1"
    "X9"             = "9"
    "X10"            = "10"
    "X11"            = "11"
    "X12"            = "12"
    "X3"             = "3"
    "X4"             = "4"
    "X5"             = "5"
    "X6"             = "6"
    "X7"             = "7"
    "X8"             = "8"
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    # The labels that need renaming only ever appear in the header row
    # (the first row of each sheet's used range); data rows below are
    # plain numbers and are left untouched.
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($startRow, $startCol + $c)
        $old = $cell.Value()
        if ($old -ne $null -and $map.ContainsKey([string]$old)) {
            $new = $map[[string]$old]

            # Some replacement values look numeric ("9", "10", ...).
            # Force the cell to stay text so it keeps being written as
            # a shared string instead of turning into a numeric cell.
            $looksNumeric = $new -match '^[0-9]+$'
            if ($looksNumeric) {
                $cell.NumberFormat = "@"
            }

            $cell.Value = $new

            if ($looksNumeric) {
                # Drop the temporary text format again so no stray
                # number-format/style gets left behind on the cell.
                $cell.Style = "Normal"
            }
        }
    }
}

Write-Output "Renamed header labels on $($wb.Worksheets.Count) worksheet(s)."
